# Working-hours log: add the 2014-03-29 entries that were missing, pushing
# the "empty separator / sum [min] / sum [h] / sum [working weeks]" block
# down by two rows (it now lives at rows 102-105 instead of 100-103).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the old separator row (100), so the
# trailing summary rows (previously 100-103) shift down to 102-105 and
# inherit the D/E/F/G number formats from the row above, exactly like
# typing new rows in the UI would.
$ws.Rows("100:101").Insert()

# 2014-03-29, 13:00 - 13:30
$ws.Range("A100").Value = 2014
$ws.Range("B100").Value = 3
$ws.Range("C100").Value = 29
$ws.Range("D100").Value = 0.54166666666666663
$ws.Range("E100").Value = 0.5625
$ws.Range("F100").Formula = "=(E100-D100)*24*60"
$ws.Range("G100").Formula = "=F100/60"

# 2014-03-29, 16:00 - 19:00
$ws.Range("A101").Value = 2014
$ws.Range("B101").Value = 3
$ws.Range("C101").Value = 29
$ws.Range("D101").Value = 0.66666666666666663
$ws.Range("E101").Value = 0.79166666666666663
$ws.Range("F101").Formula = "=(E101-D101)*24*60"
$ws.Range("G101").Formula = "=F101/60"

# The sum formulas in the (now shifted) summary rows already point at the
# right ranges because Excel adjusted them automatically on row insert;
# only the view selection needs to move to where the cursor ended up.
$ws.Range("A102").Select()
